# Update "想去人数" (interested-people count, column F) values on the
# 展览 (sheet1) and 全部类型 (sheet4) sheets, and the one value that also
# changed on 本地生活 (sheet3), to reflect a refreshed scrape of the site.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 sheet (sheet1)
$ws1.Range("F2").Value = 538
$ws1.Range("F3").Value = 927
$ws1.Range("F30").Value = 308
$ws1.Range("F36").Value = 3967
$ws1.Range("F37").Value = 58

# 本地生活 sheet (sheet3)
$ws3.Range("F5").Value = 1666

# 全部类型 sheet (sheet4) - mirrors the rows above
$ws4.Range("F4").Value = 1666
$ws4.Range("F7").Value = 538
$ws4.Range("F8").Value = 927
$ws4.Range("F42").Value = 308
$ws4.Range("F50").Value = 3967
$ws4.Range("F51").Value = 58
